# Append 12 new alert rows (653-664) to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows, in column order A..O.
# Columns: stockname, detected_date, rowNumber, date1, row1, value1, date2, row2, value2,
#          buyORsell, slope, intercept, window_size, percentage_value, two_line_count
$rowsCsv = @"
BATAINDIA.NS|39417|65|38869|47|59.47750667254579|39295|61|59.69378497807872|Low|0.01544845039520989|58.75142950397092|3|1|2
MANAPPURAM.NS|42917|84|42644|75|90.23010662436555|42795|80|90.57228138015006|High|0.06843495115690246|85.09748528759786|3|1|2
GODREJPROP.NS|43344|103|43101|95|911.8499755859375|43221|99|920|High|2.037506103515625|718.2868957519531|3|1|2
GMRAIRPORT.NS|42370|112|42217|107|8.829098701477051|42248|108|8.829098701477051|Low|0|8.829098701477051|3|1|2
GMRAIRPORT.NS|42522|117|42217|107|8.829098701477051|42401|113|8.784051895141602|Low|-0.007507801055908203|9.632433414459229|3|1|2
GMRAIRPORT.NS|42522|117|42248|108|8.829098701477051|42401|113|8.784051895141602|Low|-0.009009361267089844|9.802109718322754|3|1|2
BATAINDIA.NS|42705|173|42217|157|581.4596786134197|42583|169|577.2650448407002|High|-0.3495528143932916|636.3394704731664|3|1|2
MANAPPURAM.NS|44348|131|44013|120|170.1301489508114|44228|127|168.705578165982|High|-0.2035101121184937|194.5513624050307|3|1|2
LICHSGFIN.NS|42461|165|42217|157|431.6130438757089|42339|161|428.005344062932|High|-0.9019249531942393|573.2152615272045|3|1|2
LICHSGFIN.NS|42522|167|42156|155|320.6515052936053|42401|163|323.2496289079498|Low|0.3247654517930698|270.3128602656794|3|1|2
GMRAIRPORT.NS|43983|165|43344|144|19.36996269226074|43525|150|19.23482322692871|High|-0.02252324422200521|22.61330986022949|3|1|2
GMRAIRPORT.NS|44013|166|43617|153|12.65804481506348|43891|162|12.70309066772461|Low|0.005005094740125868|11.89226531982422|3|1|2
"@

$lines = $rowsCsv -split "`n" | Where-Object { $_.Trim() -ne "" }

$startRow = 653
$dateNumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $lines.Count; $i++) {
    $fields = $lines[$i].Trim() -split '\|'
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $fields[0]                      # A stockname
    $ws.Cells.Item($r, 2).Value = [double]$fields[1]               # B detected_date
    $ws.Cells.Item($r, 3).Value = [double]$fields[2]               # C rowNumber
    $ws.Cells.Item($r, 4).Value = [double]$fields[3]               # D date1
    $ws.Cells.Item($r, 5).Value = [double]$fields[4]               # E row1
    $ws.Cells.Item($r, 6).Value = [double]$fields[5]               # F value1
    $ws.Cells.Item($r, 7).Value = [double]$fields[6]               # G date2
    $ws.Cells.Item($r, 8).Value = [double]$fields[7]               # H row2
    $ws.Cells.Item($r, 9).Value = [double]$fields[8]               # I value2
    $ws.Cells.Item($r, 10).Value = $fields[9]                      # J buyORsell
    $ws.Cells.Item($r, 11).Value = [double]$fields[10]             # K slope
    $ws.Cells.Item($r, 12).Value = [double]$fields[11]             # L intercept
    $ws.Cells.Item($r, 13).Value = [double]$fields[12]             # M window_size
    $ws.Cells.Item($r, 14).Value = [double]$fields[13]             # N percentage_value
    $ws.Cells.Item($r, 15).Value = [double]$fields[14]             # O two_line_count

    # Date columns B, D, G carry the custom date/time number format used
    # throughout the sheet for serial-date values.
    $ws.Cells.Item($r, 2).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($r, 4).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($r, 7).NumberFormat = $dateNumberFormat
}

